$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-10 Sunday" "2025-08-11 Monday"
Replace-Text "84×42=3528" "41×46=1886"
Replace-Text "61×77=4697" "44×95=4180"
Replace-Text "95×30=2850" "85×68=5780"
Replace-Text "40×67=2680" "61×32=1952"
Replace-Text "76×78=5928" "36×76=2736"
Replace-Text "34×57=1938" "24×20=480"
Replace-Text "33×99=3267" "51×75=3825"
Replace-Text "72×96=6912" "36×38=1368"
Replace-Text "22×55=1210" "73×69=5037"
Replace-Text "16×36=576" "71×80=5680"
Replace-Text "89×54=4806" "26×15=390"
Replace-Text "66×91=6006" "28×66=1848"
Replace-Text "38×58=2204" "95×94=8930"
Replace-Text "60×37=2220" "39×47=1833"
Replace-Text "19×90=1710" "90×40=3600"
Replace-Text "26×87=2262" "41×97=3977"
Replace-Text "57×78=4446" "57×26=1482"
Replace-Text "77×54=4158" "26×30=780"
Replace-Text "17×92=1564" "72×70=5040"
Replace-Text "32×41=1312" "73×11=803"
Replace-Text "29×96=2784" "37×22=814"
Replace-Text "15×45=675" "20×55=1100"
Replace-Text "17×83=1411" "95×54=5130"
Replace-Text "35×74=2590" "43×94=4042"
Replace-Text "53×90=4770" "13×51=663"
